$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 23 (label "pc") entirely; remaining rows below shift up by one.
$ws.Rows.Item(23).Delete()

# Append a new row of data at the new last row (48) with label "zy_r".
$ws.Range("A48").Value = "zy_r"
$ws.Range("B48").Value = -54.28
$ws.Range("C48").Value = 2.21
$ws.Range("D48").Value = 74.33

# Match the author's final on-screen selection.
$ws.Range("D48").Select() | Out-Null
